$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AZ (52) currently holds the "Mean" header/values.
# We are inserting a new "Run 50" data column before the Mean column,
# so the old Mean column (AZ) shifts to the new column BA.

# 1) Clone the header cell's formatting (bold font, borders, alignment)
#    from AZ1 into the new BA1 header cell, then set its text.
$ws.Range("AZ1").Copy($ws.Range("BA1"))
$ws.Range("BA1").Value = "Mean"

# 2) Re-point AZ1's header text to the new run label.
$ws.Range("AZ1").Value = "Run 50"

# 3) Move the existing Mean values (previously in AZ) into the new BA
#    column, and fill AZ with the new Run 50 results.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 53).Value = 16.43680692
    $ws.Cells.Item($r, 52).Value = 20.23187132
}
